# Update "想去人数" (F column) counts on the 展览 sheet and the 全部类型 sheet.
# The 全部类型 sheet has one extra row (row 33) compared to 展览, so rows
# after it are shifted by +1 between the two sheets.

$wb = $excel.ActiveWorkbook

# Row => New F value, for rows 2-32 (identical row numbers on both sheets)
$commonUpdates = @{
    2  = 147
    3  = 354
    4  = 454
    5  = 1778
    7  = 2235
    9  = 288
    11 = 5059
    12 = 371
    14 = 313
    15 = 234
    16 = 37
    17 = 199
    18 = 390
    20 = 127
    21 = 4066
    22 = 734
    23 = 722
    24 = 33
    26 = 114
    27 = 132
    30 = 97
    31 = 593
    32 = 13
}

# Rows 34,35,36,38 on 展览 (sheet "展览") map to rows 35,36,37,39 on 全部类型.
$sheet1OnlyUpdates = @{
    34 = 1039
    35 = 6
    36 = 2640
    38 = 42
}

$sheet4OnlyUpdates = @{
    35 = 1039
    36 = 6
    37 = 2640
    39 = 42
}

$wsExpo = $wb.Worksheets.Item("展览")
foreach ($row in $commonUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
foreach ($row in $sheet1OnlyUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $sheet1OnlyUpdates[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $commonUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
foreach ($row in $sheet4OnlyUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $sheet4OnlyUpdates[$row]
}

$wb.Save()
